$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 currently holds:
#   B11 = "type: blog / width: 2 / height: 1 / ser: 154"
#   D11 = "type: blog / width: 2 / height: 1 / ser: 155"
#   I11 = "type: blog / width: 2 / height: 1 / ser: 153"
#
# The edit promotes the featured blog series: ser 153 is dropped, 154 and 155
# shift up, and a brand new "ser: 156" entry takes the lead spot (B11).
$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 155"
$ws.Range("D11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 154"
$ws.Range("B11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 156"

# Update the view state to match: scrolled/selected cell moved from F11 to I11.
$ws.Range("I11").Select()
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Application.ActiveWindow.ScrollColumn = 2
